$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'sport pants women',
    'sport pants women plus size',
    'sport resistance',
    'sport running',
    'sport running pants women',
    'sport skirts for women knee length',
    'sport skirts for women with leggings',
    'sport support',
    'sport support tape',
    'sport tape for knee',
    'sport tape kinesiology for knee',
    'sport tight',
    'sport tights',
    'sport wear',
    'sport woman pants',
    'sporting apparel',
    'sports bra cwx',
    'sports capri leggings women',
    'sports capri pants for women',
    'sports compression leggings',
    'sports compression pants',
    'sports compression shorts',
    'sports for women leggings',
    'sports gear for women',
    'sports joint',
    'sports knee support',
    'sports leg compression',
    'sports leggings for women',
    'sports leggings women capri',
    'sports pants for women',
    'sports recovery equipment',
    'sports running pants',
    'sports stocking',
    'sports tape athletic kinesiology',
    'sports tape for knee',
    'sports tape for pain',
    'sports tape kinesiology for knees',
    'sports taping techniques',
    'sports team wear',
    'sports teams fabric',
    'sports tear away pants',
    'sports tights',
    'sports wear for women gym',
    'sports wear women plus size',
    'sports women nike',
    'spring tights',
    'squad leggings',
    'squat compression shorts',
    'squat leggings women',
    'squat proof high waisted leggings',
    'squat proof leggings',
    'squat proof shorts women',
    'squat squad',
    'squeeze ball for thighs',
    'squeeze band',
    'squeeze key ring',
    'st joseph basketball',
    'stability ball chart',
    'stability ball exercise chart',
    'stability ball legs',
    'stability pro',
    'stability pro back support',
    'stability pro knee brace',
    'stability pro knee compression',
    'stability training',
    'stabilization ball',
    'stabilizing grip',
    'stabilizing stick',
    'stabilizing tape',
    'stabilyx joint support',
    'stabilyx joint support compression tight',
    'stabilyx men',
    'stabilyx tights women',
    'stabylix pants women',
    'stage clothing',
    'stage compression garment',
    'stage wear for men',
    'stand by generators',
    'standing seam',
    'starter leggings women',
    'statement pants',
    'stay hip',
    'stay on target',
    'stay up leggings',
    'stay up thigh high hose',
    'stay warm leggings',
    'steep skiing',
    'steep snowboarding',
    'still cool running speed training',
    'still i rise',
    'stocking black for women',
    'stocking plus size women',
    'stockings and tights',
    'stockings knee high',
    'stockings knee highs',
    'stockings long',
    'stockings women',
    'stomach band for after surgery',
    'stomach compression',
    'stomach compression for men'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
